# Apply the crypto-list refresh: updated prices/volume deltas for existing rows,
# plus the Hedera/WEMIXToken row swap (rows 34-35), per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "41.655.64"
    "E2" = "  +0.32%  "
    "D3" = "2.473.19"
    "E3" = "  +0.01%  "
    "D4" = "0.998"
    "E4" = "  -0.19%  "
    "D5" = "317.48"
    "E5" = "  +1.44%  "
    "D6" = "92.75"
    "E6" = "  +1.30%  "
    "E7" = "  +0.95%  "
    "E8" = "  +0.03%  "
    "E9" = "  +1.21%  "
    "D10" = "33.06"
    "E10" = "  +1.94%  "
    "D11" = "0.0853"
    "E11" = "  +8.45%  "
    "E12" = "  +0.57%  "
    "D13" = "2.856.66"
    "E13" = "  -0.01%  "
    "D14" = "6.90"
    "E14" = "  +0.66%  "
    "D15" = "15.75"
    "E15" = "  -3.16%  "
    "D16" = "2.493.99"
    "E16" = "  +1.45%  "
    "D17" = "0.791"
    "E17" = "  +2.82%  "
    "D18" = "41.616.96"
    "E18" = "  +0.26%  "
    "D19" = "6.46"
    "E19" = "  -0.33%  "
    "D20" = "0.0₃0948"
    "E20" = "  +0.56%  "
    "D21" = "71.33"
    "E21" = "  -0.58%  "
    "D22" = "11.28"
    "D23" = "239.42"
    "E23" = "  +1.51%  "
    "E24" = "  +1.16%  "
    "E25" = "  +1.74%  "
    "D26" = "0.999"
    "E26" = "  -0.13%  "
    "D27" = "24.71"
    "E27" = "  -0.22%  "
    "E28" = "  +2.54%  "
    "D29" = "9.84"
    "D30" = "36.21"
    "E30" = "  +1.83%  "
    "D31" = "160.12"
    "E31" = "  +2.62%  "
    "E32" = "  +1.60%  "
    "E33" = "  -0.06%  "
    "B34" = "WEMIXToken"
    "C34" = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
    "D34" = "2.59"
    "E34" = "  +0.55%  "
    "B35" = "Hedera"
    "C35" = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
    "D35" = "0.0768"
    "E35" = "  +1.70%  "
    "D36" = "17.32"
    "E36" = "  +0.15%  "
    "D37" = "2.92"
    "E37" = "  +1.55%  "
    "E38" = "  +2.17%  "
    "E39" = "  +1.81%  "
    "E40" = "  -0.45%  "
    "E41" = "  -1.27%  "
    "E42" = "  +3.08%  "
    "D43" = "1.984.83"
    "E43" = "  +1.38%  "
    "E44" = "  +0.61%  "
    "D45" = "18.81"
    "E45" = "  +0.22%  "
    "D46" = "2.99"
    "E46" = "  +2.19%  "
    "E47" = "  +3.69%  "
    "D48" = "2.714.79"
    "E48" = "  -0.04%  "
    "D49" = "97.29"
    "E49" = "  -0.55%  "
    "D50" = "74.07"
    "E50" = "  +2.85%  "
    "D51" = "67.11"
    "E51" = "  -0.19%  "
}

foreach ($cellRef in $updates.Keys) {
    # Leading apostrophe forces text interpretation so numeric-looking strings
    # (e.g. "41.655.64", "0.998", "0.0₃0948") are kept as text, matching the
    # original inline-string cell type instead of being coerced to a number.
    $ws.Range($cellRef).Value = "'" + $updates[$cellRef]
    # Clear the quote-prefix formatting flag so no stray number format/style
    # is introduced on the cell (keeps it on the default/unstyled xf).
    $ws.Range($cellRef).Style = "Normal"
}
